$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Update the letter date: "September 19, 2025" -> "September 21, 2025"
#    (Paragraph.Range.Text includes the trailing paragraph mark, so we
#    match with -like rather than -eq.)
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "September 19, 2025*") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# ------------------------------------------------------------------
# 2) Split the mailing-address paragraph ("909 Story Road, San Jose CA
#    95122") into two paragraphs: "909 Story Road" and "San Jose, CA
#    95122". Only the standalone address paragraph (not the one inside
#    the summary table) should be split, so locate it by scanning the
#    document paragraphs for an exact text match.
# ------------------------------------------------------------------
$addrIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "909 Story Road, San Jose CA 95122*") {
        $addrIndex = $i
        break
    }
}

if ($addrIndex -ne -1) {
    $addrPara = $d.Paragraphs.Item($addrIndex)
    # Shrink the existing run down to the street line (keeps the run's
    # original formatting, including xml:space="preserve").
    $addrPara.Range.Text = "909 Story Road"
    # Insert a brand-new paragraph right after it, using the same
    # paragraph formatting, and fill it with the city/state/zip line.
    $null = $addrPara.Range.InsertParagraphAfter()
    $cityPara = $d.Paragraphs.Item($addrIndex + 1)
    $cityPara.Range.Text = "San Jose, CA 95122"
    $cityPara.Range.Font.Name = "Arial"
    $cityPara.Range.Font.Size = 11
}

# ------------------------------------------------------------------
# 3) Remove the now-superfluous blank "No Spacing" paragraph that sits
#    directly after "Board of Directors".
# ------------------------------------------------------------------
$bodIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Board of Directors*") {
        $bodIndex = $i
        break
    }
}

if ($bodIndex -ne -1) {
    $nextPara = $d.Paragraphs.Item($bodIndex + 1)
    if ($nextPara.Range.Text -eq [char]13) {
        $nextPara.Range.Delete()
    }
}
